$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

# A new contract-note transaction (NSE buy on 2026-02-13) needs to be
# inserted above the current top data row, pushing the existing
# transactions (currently rows 5-7) down to rows 6-8.
$ws.Rows.Item(5).Insert()

# Excel's row Insert() copies the formatting of the row above (the bold
# header row), so strip that back off before writing the new data - the
# data rows in this sheet carry no explicit formatting except the date
# column.
$ws.Rows.Item(5).ClearFormats()

# Column A keeps the same custom date/time number format used by the
# other transaction rows.
$ws.Range("A5").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("A5").Value = 46066
$ws.Range("B5").Value = "NSE"
$ws.Range("C5").Value = "Buy"
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 952.35
$ws.Range("F5").Value = 4795.29
$ws.Range("G5").Value = "CN#252611910666"
$ws.Range("H5").Value = 4.81
$ws.Range("I5").Value = 28.732
$ws.Range("J5").Formula = '=Index!$C$2'

# The Insert() above also stamped the header's formatting (and therefore
# empty styled cells) across the rest of row 5's width (K:N, W:AB) even
# though those columns have no data on transaction rows. Clear them fully
# so they don't linger as empty cells in the saved file.
$ws.Range("K5:N5").Clear()
$ws.Range("W5:AB5").Clear()
